$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header in column H, to the right of the existing "sum" header (G1).
$ws.Range("H1").Value = "Save"

# Copy G1's formatting (bold font, border, centered alignment) onto H1 so the
# new header matches the look of the existing header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new "Save" column (rows 2-9) with 0 for every existing record.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
